# Constant folding pass in semantic analyzer (#23)
#
# Adds two new SemAnalyzer error descriptions (rows 15 & 16 of the Errors
# sheet) and moves the active selection to D17 (the cell right after the
# newly-filled-in rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Errors")

$ws.Range("D15").Value = "Integer overflow (max 999 min -999)"
$ws.Range("D16").Value = "Div 0/Mod0"

$ws.Activate()
$ws.Range("D17").Select()
